# Apply the "evolution diff date format YYYY, YYYYtQ and YYYY/MM" edit.
#
# This updates a handful of values in the dataset table on Sheet1:
#   - row 7  (eco_energies)   : start_date "2019" -> "2019/07" (YYYY/MM)
#                                end_date   "2024" -> "2024t2"  (YYYYtQ)
#   - row 8  (eau_potable)    : nb_row 80000 -> 34444
#   - row 11 (transport_pub)  : nb_row 90000 -> 999990
#   - row 18 (conso_energie)  : start_date "2015" -> "2011"
#   - row 19 (revenus_menages): end_date   "2024" -> "2027"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# eco_energies (row 7): evolve start_date/end_date formats
$ws.Range("L7").Value = "2019/07"
$ws.Range("M7").Value = "2024t2"

# eau_potable (row 8): updated row count
$ws.Range("H8").Value = 34444

# transport_pub (row 11): updated row count
$ws.Range("H11").Value = 999990

# conso_energie (row 18): updated start_date
$ws.Range("L18").Value = "2011"

# revenus_menages (row 19): updated end_date
$ws.Range("M19").Value = "2027"

# Match the saved selection from the authored edit (last touched cell).
$ws.Range("L18").Select()
